$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 564, pushing the existing rows 564-656 down to 565-657.
$ws.Rows(564).Insert()

# Populate the newly inserted row 564 with the new weekly record.
$ws.Cells.Item(564, 1).Value  = 3
$ws.Cells.Item(564, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(564, 3).Value  = "Coquimbo"
$ws.Cells.Item(564, 4).Value  = 45218
$ws.Cells.Item(564, 5).Value  = 5
$ws.Cells.Item(564, 6).Value  = 100112043
$ws.Cells.Item(564, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(564, 8).Value  = "Sin especificar"
$ws.Cells.Item(564, 9).Value  = "Primera"
$ws.Cells.Item(564, 10).Value = 68
$ws.Cells.Item(564, 11).Value = 12000
$ws.Cells.Item(564, 12).Value = 12000
$ws.Cells.Item(564, 13).Value = 12000
$ws.Cells.Item(564, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(564, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(564, 16).Value = 200
$ws.Cells.Item(564, 17).Value = 60
$ws.Cells.Item(564, 18).Value = "Hortaliza"
